# "Generate Report for Archive"
#
# The localization status report is being regenerated: the outstanding
# "Ready for handoff" status is updated to "In Translation" everywhere it
# appears (Overview summary columns + each per-locale sheet's Status
# column), and the Status column on each sheet is re-sized to fit the new,
# shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells (columns E and F) -------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"

# --- Per-locale sheets: Status column (column C) -------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value2 = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value2 = "In Translation"

# --- Resize the status columns to fit the new text ------------------------
# ColumnWidth is expressed in "characters of the Normal style font" and
# Excel snaps it to whole-pixel boundaries, so we hand it the value that
# lands on the narrower width the shorter "In Translation" label now needs.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
